# Auto-generated market-data refresh for Asura_Profits workbook
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) per sheet
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 31996.188
$ws.Range("I74").Value = 46456.094
$ws.Range("J74").Value = 4390.909
$ws.Range("K74").Value = 46456.094
$ws.Range("L74").Value = 4390.909
$ws.Range("M74").Value = -45520.094
$ws.Range("N74").Value = -6262.909
$ws.Range("H77").Value = 31996.188
$ws.Range("I77").Value = 46456.094
$ws.Range("J77").Value = 4390.909
$ws.Range("K77").Value = 232280.47
$ws.Range("L77").Value = 21954.545
$ws.Range("M77").Value = -227600.47
$ws.Range("N77").Value = -31314.545
$ws.Range("H125").Value = 6120
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 8866.666999999999
$ws.Range("K125").Value = 18000
$ws.Range("L125").Value = 79800.003
$ws.Range("M125").Value = -15540
$ws.Range("N125").Value = -84720.003
$ws.Range("H135").Value = 559.1951
$ws.Range("I135").Value = 582.9
$ws.Range("J135").Value = 494.54544
$ws.Range("K135").Value = 5246.099999999999
$ws.Range("L135").Value = 4450.90896
$ws.Range("M135").Value = -2711.099999999999
$ws.Range("N135").Value = -9520.908960000001
$ws.Range("H137").Value = 2040.5625
$ws.Range("I137").Value = 1318.9231
$ws.Range("J137").Value = 2893.4092
$ws.Range("K137").Value = 3956.7693
$ws.Range("L137").Value = 8680.2276
$ws.Range("M137").Value = -1406.7693
$ws.Range("N137").Value = -13780.2276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1160.0518
$ws.Range("I74").Value = 1183.907
$ws.Range("J74").Value = 1091.6666
$ws.Range("K74").Value = 1183.907
$ws.Range("L74").Value = 1091.6666
$ws.Range("M74").Value = -309.9069999999999
$ws.Range("N74").Value = -2839.6666
$ws.Range("H77").Value = 1160.0518
$ws.Range("I77").Value = 1183.907
$ws.Range("J77").Value = 1091.6666
$ws.Range("K77").Value = 5919.535
$ws.Range("L77").Value = 5458.333000000001
$ws.Range("M77").Value = -1551.535
$ws.Range("N77").Value = -14194.333
$ws.Range("H123").Value = 25427.715
$ws.Range("J123").Value = 25427.715
$ws.Range("L123").Value = 25427.715
$ws.Range("N123").Value = -35227.715
$ws.Range("H131").Value = 44285.6
$ws.Range("J131").Value = 44285.6
$ws.Range("L131").Value = 44285.6
$ws.Range("N131").Value = -54365.6
$ws.Range("H132").Value = 4907.341
$ws.Range("I132").Value = 5448.793
$ws.Range("J132").Value = 3860.5334
$ws.Range("K132").Value = 16346.379
$ws.Range("L132").Value = 11581.6002
$ws.Range("M132").Value = -13816.379
$ws.Range("N132").Value = -16641.6002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1876.7433
$ws.Range("I31").Value = 2517.7097
$ws.Range("J31").Value = 1414.6511
$ws.Range("K31").Value = 2517.7097
$ws.Range("L31").Value = 1414.6511
$ws.Range("M31").Value = -2222.7097
$ws.Range("N31").Value = -2004.6511
$ws.Range("H34").Value = 1876.7433
$ws.Range("I34").Value = 2517.7097
$ws.Range("J34").Value = 1414.6511
$ws.Range("K34").Value = 2517.7097
$ws.Range("L34").Value = 1414.6511
$ws.Range("M34").Value = -2315.7097
$ws.Range("N34").Value = -1818.6511
$ws.Range("H134").Value = 1529.9615
$ws.Range("I134").Value = 1295
$ws.Range("J134").Value = 4349.5
$ws.Range("K134").Value = 3885
$ws.Range("L134").Value = 13048.5
$ws.Range("M134").Value = -1350
$ws.Range("N134").Value = -18118.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H34").Value = 980.6667
$ws.Range("I34").Value = 171
$ws.Range("J34").Value = 2600
$ws.Range("K34").Value = 513
$ws.Range("L34").Value = 7800
$ws.Range("M34").Value = -429
$ws.Range("N34").Value = -7968
$ws.Range("H63").Value = 236602.67
$ws.Range("J63").Value = 5000
$ws.Range("L63").Value = 15000
$ws.Range("N63").Value = -16498
$ws.Range("H66").Value = 236602.67
$ws.Range("J66").Value = 5000
$ws.Range("L66").Value = 45000
$ws.Range("N66").Value = -52488
$ws.Range("H76").Value = 10071
$ws.Range("I76").Value = 213
$ws.Range("J76").Value = 15000
$ws.Range("K76").Value = 639
$ws.Range("L76").Value = 45000
$ws.Range("M76").Value = -256
$ws.Range("N76").Value = -45766
$ws.Range("H79").Value = 10071
$ws.Range("I79").Value = 213
$ws.Range("J79").Value = 15000
$ws.Range("K79").Value = 639
$ws.Range("L79").Value = 45000
$ws.Range("M79").Value = 687
$ws.Range("N79").Value = -47652
$ws.Range("H107").Value = 1235.3125
$ws.Range("I107").Value = 1190.9231
$ws.Range("J107").Value = 1304.56
$ws.Range("K107").Value = 3572.7693
$ws.Range("L107").Value = 3913.68
$ws.Range("M107").Value = -1652.7693
$ws.Range("N107").Value = -7753.68
$ws.Range("H122").Value = 670.06665
$ws.Range("I122").Value = 580.9231
$ws.Range("J122").Value = 1249.5
$ws.Range("K122").Value = 5228.3079
$ws.Range("L122").Value = 11245.5
$ws.Range("M122").Value = -2778.3079
$ws.Range("N122").Value = -16145.5
$ws.Range("H131").Value = 3956.0571
$ws.Range("I131").Value = 625.9091
$ws.Range("J131").Value = 5482.375
$ws.Range("K131").Value = 1877.7273
$ws.Range("L131").Value = 16447.125
$ws.Range("M131").Value = 3162.2727
$ws.Range("N131").Value = -26527.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 8609.625
$ws.Range("J123").Value = 8609.625
$ws.Range("L123").Value = 8609.625
$ws.Range("N123").Value = -13509.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 43333.332
$ws.Range("J109").Value = 43333.332
$ws.Range("L109").Value = 43333.332
$ws.Range("N109").Value = -46107.332
$ws.Range("H136").Value = 1698.5193
$ws.Range("I136").Value = 1262.1562
$ws.Range("J136").Value = 2396.7
$ws.Range("K136").Value = 3786.4686
$ws.Range("L136").Value = 7190.099999999999
$ws.Range("M136").Value = -1236.4686
$ws.Range("N136").Value = -12290.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 104352.37
$ws.Range("I81").Value = 141109.25
$ws.Range("J81").Value = 6334
$ws.Range("K81").Value = 282218.5
$ws.Range("L81").Value = 12668
$ws.Range("M81").Value = -281157.5
$ws.Range("N81").Value = -14790
$ws.Range("H84").Value = 104352.37
$ws.Range("I84").Value = 141109.25
$ws.Range("J84").Value = 6334
$ws.Range("K84").Value = 1411092.5
$ws.Range("L84").Value = 63340
$ws.Range("M84").Value = -1405788.5
$ws.Range("N84").Value = -73948
$ws.Range("H107").Value = 366.3846
$ws.Range("I107").Value = 269.45456
$ws.Range("J107").Value = 899.5
$ws.Range("K107").Value = 808.36368
$ws.Range("L107").Value = 2698.5
$ws.Range("M107").Value = 1111.63632
$ws.Range("N107").Value = -6538.5
$ws.Range("H108").Value = 45313
$ws.Range("J108").Value = 45313
$ws.Range("L108").Value = 45313
$ws.Range("N108").Value = -52993
$ws.Range("H119").Value = 500698
$ws.Range("J119").Value = 500698
$ws.Range("L119").Value = 500698
$ws.Range("N119").Value = -510374
$ws.Range("H124").Value = 500429
$ws.Range("J124").Value = 500429
$ws.Range("L124").Value = 500429
$ws.Range("N124").Value = -510249
$ws.Range("H126").Value = 7554.4546
$ws.Range("I126").Value = 8805.666999999999
$ws.Range("K126").Value = 26417.001
$ws.Range("M126").Value = -23947.001
$ws.Range("H129").Value = 35428
$ws.Range("J129").Value = 35428
$ws.Range("L129").Value = 35428
$ws.Range("N129").Value = -45428
$ws.Range("H131").Value = 54614.168
$ws.Range("J131").Value = 54614.168
$ws.Range("L131").Value = 54614.168
$ws.Range("N131").Value = -64694.168
$ws.Range("H132").Value = 1935.5116
$ws.Range("I132").Value = 1822.8966
$ws.Range("J132").Value = 2168.7856
$ws.Range("K132").Value = 5468.6898
$ws.Range("L132").Value = 6506.3568
$ws.Range("M132").Value = -2938.6898
$ws.Range("N132").Value = -11566.3568
